# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with freshly scraped figures, mirroring the automated
# "Updated cryptos list ... with GitHub Actions" commit.

function Looks-Numeric($val) {
    # Matches plain decimal numbers such as "576.43" or "-0.12" that Excel
    # would otherwise silently reinterpret as a numeric value when assigned
    # through Range.Value (losing the original text formatting / exact
    # string representation used throughout this sheet).
    return $val -match '^[+-]?[0-9]+(\.[0-9]+)?$'
}

function Set-CellText($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    if (Looks-Numeric $val) {
        # Force a text number format first so the assigned string is kept
        # verbatim as text instead of being auto-converted into a number.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "66.854.10"
Set-CellText $ws "E2" "  -0.45%  "

Set-CellText $ws "D3" "3.103.86"
Set-CellText $ws "E3" "  -0.12%  "

Set-CellText $ws "E4" "  -0.10%  "

Set-CellText $ws "D5" "576.43"
Set-CellText $ws "E5" "  -0.48%  "

Set-CellText $ws "D6" "177.53"
Set-CellText $ws "E6" "  +2.22%  "

Set-CellText $ws "E7" "  -0.03%  "

Set-CellText $ws "D8" "3.102.86"
Set-CellText $ws "E8" "  +0.04%  "

Set-CellText $ws "E9" "  -1.13%  "

Set-CellText $ws "D10" "6.34"
Set-CellText $ws "E10" "  -3.08%  "

Set-CellText $ws "E11" "  -0.90%  "

Set-CellText $ws "E12" "  -1.85%  "

Set-CellText $ws "E13" "  -2.99%  "

Set-CellText $ws "E14" "  -1.74%  "

Set-CellText $ws "E15" "  -0.54%  "

Set-CellText $ws "D16" "3.619.83"
Set-CellText $ws "E16" "  -0.10%  "

Set-CellText $ws "D17" "66.850.98"
Set-CellText $ws "E17" "  -0.46%  "

Set-CellText $ws "E18" "  -0.64%  "

Set-CellText $ws "D19" "3.101.77"
Set-CellText $ws "E19" "  -0.37%  "

Set-CellText $ws "D20" "16.69"
Set-CellText $ws "E20" "  +0.77%  "

Set-CellText $ws "D21" "480.43"
Set-CellText $ws "E21" "  -2.18%  "

Set-CellText $ws "D22" "7.82"
Set-CellText $ws "E22" "  -0.19%  "

Set-CellText $ws "D23" "0.691"
Set-CellText $ws "E23" "  -1.53%  "

Set-CellText $ws "D24" "83.66"
Set-CellText $ws "E24" "  -0.33%  "

Set-CellText $ws "E25" "  -3.50%  "

Set-CellText $ws "D26" "2.24"
Set-CellText $ws "E26" "  -1.89%  "

Set-CellText $ws "E27" "  -4.62%  "

Set-CellText $ws "E28" "  +0.06%  "

Set-CellText $ws "D29" "7.91"
Set-CellText $ws "E29" "  -0.04%  "

Set-CellText $ws "E30" "  -2.92%  "

Set-CellText $ws "D31" "2.60"
Set-CellText $ws "E31" "  -2.17%  "

Set-CellText $ws "D32" "27.96"
Set-CellText $ws "E32" "  -1.25%  "

Set-CellText $ws "D33" "0.111"
Set-CellText $ws "E33" "  -1.78%  "

Set-CellText $ws "D34" "0.0₃0936"
Set-CellText $ws "E34" "  +0.11%  "

Set-CellText $ws "E35" "  -0.14%  "

Set-CellText $ws "D36" "48.26"
Set-CellText $ws "E36" "  +2.38%  "

Set-CellText $ws "E37" "  -4.16%  "

Set-CellText $ws "D38" "0.939"
Set-CellText $ws "E38" "  -3.20%  "

Set-CellText $ws "D39" "0.312"
Set-CellText $ws "E39" "  +1.30%  "

Set-CellText $ws "E40" "  -2.20%  "

Set-CellText $ws "E41" "  -1.78%  "

Set-CellText $ws "E42" "  -0.34%  "

Set-CellText $ws "E43" "  -1.39%  "

Set-CellText $ws "E44" "  +4.67%  "

Set-CellText $ws "D45" "2.801.24"
Set-CellText $ws "E45" "  +0.10%  "

Set-CellText $ws "D46" "374.73"
Set-CellText $ws "E46" "  -3.57%  "

Set-CellText $ws "E47" "  -2.05%  "

Set-CellText $ws "D48" "135.28"
Set-CellText $ws "E48" "  +0.25%  "

Set-CellText $ws "E49" "  +0.00%  "

Set-CellText $ws "D50" "25.52"
Set-CellText $ws "E50" "  +1.65%  "

Set-CellText $ws "E51" "  +2.00%  "
